$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$r = $ws.Range("C24")
$r.Interior.Color = 255
$r.Borders.LineStyle = 1
Write-Host "ok"
